# Update market-price / profit figures in the Leve-profit sheets.
# Values were refreshed by the scheduled market-data runner; this script
# rewrites the affected cells (currentAveragePrice* / LevePrice* / LeveProfit*)
# to their new snapshot values for each (sheet, row) pair below.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 1877.5
$ws.Range("I19").Value = 829.8
$ws.Range("K19").Value = 829.8
$ws.Range("M19").Value = -654.8

# ALC row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 2964.8333
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2964.8333
$ws.Range("K32").Value = 0
$ws.Range("L32").ClearContents()
$ws.Range("M32").Value = 2964.8333
$ws.Range("N32").Value = -3616.8333

# ALC row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 205.29411
$ws.Range("I33").Value = 124.375
$ws.Range("K33").Value = 124.375
$ws.Range("M33").Value = 104.625

# ALC row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 1414
$ws.Range("I43").Value = 1333
$ws.Range("J43").Value = 1444.375
$ws.Range("K43").Value = 1333
$ws.Range("L43").Value = 1444.375
$ws.Range("M43").Value = -1264
$ws.Range("N43").Value = -1582.375

# ALC row 52 (Leve Item ID 4567)
$ws.Range("H52").Value = 800
$ws.Range("I52").Value = 800
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 2400
$ws.Range("L52").ClearContents()
$ws.Range("M52").Value = -2240
$ws.Range("N52").Value = 0

# ALC row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 35714850
$ws.Range("I92").Value = 50000388
$ws.Range("J92").Value = 997.5
$ws.Range("K92").Value = 50000388
$ws.Range("L92").Value = 997.5
$ws.Range("M92").Value = -49999140
$ws.Range("N92").Value = -3493.5

# ALC row 97 (Leve Item ID 19885)
$ws.Range("H97").Value = 1136.6
$ws.Range("J97").Value = 1198.75
$ws.Range("L97").Value = 3596.25
$ws.Range("N97").Value = -4588.25

# ALC row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2453.603
$ws.Range("J138").Value = 2203.0264
$ws.Range("L138").Value = 6609.0792
$ws.Range("N138").Value = -16889.0792

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 397316.56
$ws.Range("I2").Value = 556032
$ws.Range("K2").Value = 556032
$ws.Range("M2").Value = -555919

# ARM row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 3803.1733
$ws.Range("I32").Value = 3160.7
$ws.Range("J32").Value = 12797.8
$ws.Range("K32").Value = 3160.7
$ws.Range("L32").Value = 12797.8
$ws.Range("M32").Value = -2873.7
$ws.Range("N32").Value = -13371.8

# ARM row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 779.58826
$ws.Range("I74").Value = 543.5484
$ws.Range("K74").Value = 543.5484
$ws.Range("M74").Value = 330.4516

# ARM row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 779.58826
$ws.Range("I77").Value = 543.5484
$ws.Range("K77").Value = 2717.742
$ws.Range("M77").Value = 1650.258

# ARM row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 397316.56
$ws.Range("I116").Value = 556032
$ws.Range("K116").Value = 556032
$ws.Range("M116").Value = -553738

# ARM row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 3215.8572
$ws.Range("I122").Value = 1106
$ws.Range("K122").Value = 3318
$ws.Range("M122").Value = -868

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 397316.56
$ws.Range("I3").Value = 556032
$ws.Range("K3").Value = 556032
$ws.Range("M3").Value = -555918

# BSM row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 2263.6
$ws.Range("I20").Value = 1936
$ws.Range("K20").Value = 1936
$ws.Range("M20").Value = -1689

# BSM row 64 (Leve Item ID 14184)
$ws.Range("H64").Value = 548
$ws.Range("I64").Value = 548
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 548
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -323

# BSM row 67 (Leve Item ID 14184)
$ws.Range("H67").Value = 548
$ws.Range("I67").Value = 548
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 548
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = 232

# BSM row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 1672.7142
$ws.Range("I99").Value = 1335
$ws.Range("J99").Value = 1926
$ws.Range("K99").Value = 1335
$ws.Range("L99").Value = 1926
$ws.Range("M99").Value = 163
$ws.Range("N99").Value = -4922

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2761.889
$ws.Range("J31").Value = 2979.5715
$ws.Range("L31").Value = 2979.5715
$ws.Range("N31").Value = -3569.5715

# CRP row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2761.889
$ws.Range("J34").Value = 2979.5715
$ws.Range("L34").Value = 2979.5715
$ws.Range("N34").Value = -3383.5715

# CRP row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 1280887
$ws.Range("I58").Value = 1978259.4
$ws.Range("K58").Value = 1978259.4
$ws.Range("M58").Value = -1978056.4

# CRP row 87 (Leve Item ID 11929)
$ws.Range("H87").Value = 39499.5
$ws.Range("J87").Value = 39499.5
$ws.Range("L87").Value = 39499.5
$ws.Range("N87").Value = -41871.5

# CRP row 90 (Leve Item ID 11929)
$ws.Range("H90").Value = 39499.5
$ws.Range("J90").Value = 39499.5
$ws.Range("L90").Value = 118498.5
$ws.Range("N90").Value = -130354.5

# CRP row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 2723
$ws.Range("I99").Value = 2574
$ws.Range("J99").Value = 2971.3333
$ws.Range("K99").Value = 2574
$ws.Range("L99").Value = 2971.3333
$ws.Range("M99").Value = -1076
$ws.Range("N99").Value = -5967.3333

# CRP row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 2723
$ws.Range("I126").Value = 2574
$ws.Range("J126").Value = 2971.3333
$ws.Range("K126").Value = 7722
$ws.Range("L126").Value = 8913.999899999999
$ws.Range("M126").Value = -5252
$ws.Range("N126").Value = -13853.9999

# CRP row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2755.2856
$ws.Range("I132").Value = 1525.6666
$ws.Range("K132").Value = 4576.9998
$ws.Range("M132").Value = -2046.9998

# CRP row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 1280887
$ws.Range("I136").Value = 1978259.4
$ws.Range("K136").Value = 5934778.199999999
$ws.Range("M136").Value = -5932228.199999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 6 (Leve Item ID 4639)
$ws.Range("H6").Value = 269.6
$ws.Range("I6").Value = 269.6
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 808.8000000000001
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -695.8000000000001

# CUL row 11 (Leve Item ID 4745)
$ws.Range("H11").Value = 687
$ws.Range("I11").Value = 635
$ws.Range("K11").Value = 1905
$ws.Range("M11").Value = -1765

# CUL row 104 (Leve Item ID 19807)
$ws.Range("H104").Value = 4749.0835
$ws.Range("J104").Value = 4749.0835
$ws.Range("L104").Value = 14247.2505
$ws.Range("N104").Value = -19489.2505

# CUL row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 25687.785
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 26620.666
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 79861.99800000001
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -89941.99800000001

$ws = $wb.Worksheets.Item("GSM")
# GSM row 15 (Leve Item ID 12018)
$ws.Range("H15").Value = 36944
$ws.Range("J15").Value = 35000
$ws.Range("L15").Value = 35000
$ws.Range("N15").Value = -35576

# GSM row 47 (Leve Item ID 4343)
$ws.Range("H47").Value = 27010.334
$ws.Range("J47").Value = 27010.334
$ws.Range("L47").Value = 27010.334
$ws.Range("N47").Value = -28146.334

# GSM row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 4844.2856
$ws.Range("J70").Value = 5127.5
$ws.Range("L70").Value = 5127.5
$ws.Range("N70").Value = -5667.5

# GSM row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 4844.2856
$ws.Range("J73").Value = 5127.5
$ws.Range("L73").Value = 5127.5
$ws.Range("N73").Value = -6999.5

# GSM row 81 (Leve Item ID 12018)
$ws.Range("H81").Value = 36944
$ws.Range("J81").Value = 35000
$ws.Range("L81").Value = 35000
$ws.Range("N81").Value = -36996

# GSM row 84 (Leve Item ID 12018)
$ws.Range("H84").Value = 36944
$ws.Range("J84").Value = 35000
$ws.Range("L84").Value = 105000
$ws.Range("N84").Value = -114984

# GSM row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 3048.5715
$ws.Range("I102").Value = 2799
$ws.Range("J102").Value = 3497.8
$ws.Range("K102").Value = 2799
$ws.Range("L102").Value = 3497.8
$ws.Range("M102").Value = -1177
$ws.Range("N102").Value = -6741.8

$ws = $wb.Worksheets.Item("LTW")
# LTW row 88 (Leve Item ID 10961)
$ws.Range("H88").Value = 7594.5
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 5189
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 5189
$ws.Range("M88").Value = -9572
$ws.Range("N88").Value = -6045

# LTW row 91 (Leve Item ID 10961)
$ws.Range("H91").Value = 7594.5
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 5189
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 5189
$ws.Range("M91").Value = -8518
$ws.Range("N91").Value = -8153

# LTW row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 1771.625
$ws.Range("I136").Value = 1454.1578
$ws.Range("J136").Value = 2978
$ws.Range("K136").Value = 4362.4734
$ws.Range("L136").Value = 8934
$ws.Range("M136").Value = -1812.4734
$ws.Range("N136").Value = -14034

$ws = $wb.Worksheets.Item("WVR")
# WVR row 70 (Leve Item ID 11979)
$ws.Range("H70").Value = 37110
$ws.Range("J70").Value = 37110
$ws.Range("L70").Value = 37110
$ws.Range("N70").Value = -37740

# WVR row 73 (Leve Item ID 11979)
$ws.Range("H73").Value = 37110
$ws.Range("J73").Value = 37110
$ws.Range("L73").Value = 37110
$ws.Range("N73").Value = -39294

# WVR row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1828.0834
$ws.Range("I132").Value = 1196.9
$ws.Range("K132").Value = 3590.7
$ws.Range("M132").Value = -1060.7

# WVR row 133 (Leve Item ID 41869)
$ws.Range("H133").Value = 64482.145
$ws.Range("J133").Value = 64789
$ws.Range("L133").Value = 64789
$ws.Range("N133").Value = -74909

# WVR row 139 (Leve Item ID 43312)
$ws.Range("H139").Value = 71813.164
$ws.Range("J139").Value = 71813.164
$ws.Range("L139").Value = 71813.164
$ws.Range("N139").Value = -82093.164
